# The sheet holds a daily price table for "Cilantro" at Terminal La Palmera
# de La Serena. A new, more recent observation (11/3/2022) is inserted as
# the new row 59, pushing all subsequent rows (old 59..177) down by one
# (new 60..178), growing the used range from A1:R177 to A1:R178.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 59; Excel shifts rows 59:177 down to 60:178 and the
# new row inherits formatting (incl. the date style) from the row above.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new observation.
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = "11/3/2022"
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100112040
$ws.Cells.Item(59, 7).Value = "Cilantro"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 2600
$ws.Cells.Item(59, 11).Value = 1500
$ws.Cells.Item(59, 12).Value = 2000
$ws.Cells.Item(59, 13).Value = 1750
$ws.Cells.Item(59, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(59, 16).Value = 1167
$ws.Cells.Item(59, 17).Value = 1.5
$ws.Cells.Item(59, 18).Value = "Hortaliza"
